$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price"-column (D) cell while forcing the
# Text type (prices like "314.50"/"1.00" must stay literal text, matching
# the source data, instead of Excel auto-coercing them to numbers). The
# NumberFormat flip is reverted immediately via ClearFormats so the cell
# keeps its original (default) style.
function Set-PriceText($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Rows 11/12: Chainlink and Dogecoin swapped positions in the ranking ---
$ws.Range("B11").Value = "Chainlink"
$ws.Range("C11").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-PriceText "D11" "19.24"
$ws.Range("E11").Value = "  +4.91%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-PriceText "D12" "0.0799"
$ws.Range("E12").Value = "  +1.95%  "

# --- Refreshed price (D) / 1h-volume-change (E) figures ---
Set-PriceText "D2" "44.493.10"
$ws.Range("E2").Value = "  +3.66%  "
Set-PriceText "D3" "2.421.05"
$ws.Range("E3").Value = "  +2.51%  "
Set-PriceText "D4" "0.999"
$ws.Range("E4").Value = "  -0.03%  "
Set-PriceText "D5" "314.45"
$ws.Range("E5").Value = "  +3.85%  "
Set-PriceText "D6" "100.66"
$ws.Range("E6").Value = "  +5.55%  "
Set-PriceText "D7" "0.512"
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("E8").Value = "  -0.08%  "
Set-PriceText "D9" "0.520"
$ws.Range("E9").Value = "  +7.58%  "
Set-PriceText "D10" "35.36"
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("E14").Value = "  +3.83%  "
Set-PriceText "D15" "2.799.88"
$ws.Range("E15").Value = "  +2.57%  "
Set-PriceText "D16" "2.426.92"
$ws.Range("E16").Value = "  +2.18%  "
Set-PriceText "D17" "0.832"
$ws.Range("E17").Value = "  +5.31%  "
Set-PriceText "D18" "44.366.39"
$ws.Range("E18").Value = "  +3.47%  "
Set-PriceText "D19" "12.43"
$ws.Range("E19").Value = "  +4.94%  "
Set-PriceText "D20" "6.39"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("E21").Value = "  +3.87%  "
Set-PriceText "D22" "68.70"
$ws.Range("E22").Value = "  +1.10%  "
Set-PriceText "D23" "241.86"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("E24").Value = "  +5.63%  "
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E26").Value = "  -0.06%  "
Set-PriceText "D27" "25.14"
$ws.Range("E27").Value = "  +3.05%  "
Set-PriceText "D28" "2.28"
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("E29").Value = "  +3.17%  "
$ws.Range("E30").Value = "  +4.50%  "
Set-PriceText "D31" "48.44"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("E32").Value = "  +17.87%  "
Set-PriceText "D33" "19.25"
$ws.Range("E33").Value = "  +10.25%  "
$ws.Range("E34").Value = "  +3.56%  "
Set-PriceText "D35" "0.0773"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("E36").Value = "  +0.18%  "
Set-PriceText "D37" "1.89"
$ws.Range("E37").Value = "  +2.58%  "
Set-PriceText "D38" "4.46"
$ws.Range("E38").Value = "  +4.14%  "
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  -2.33%  "
Set-PriceText "D41" "121.19"
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("E42").Value = "  +1.57%  "
Set-PriceText "D43" "20.83"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  +4.11%  "
Set-PriceText "D45" "1.943.13"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("E47").Value = "  +8.79%  "
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("E49").Value = "  +10.29%  "
Set-PriceText "D50" "54.92"
$ws.Range("E50").Value = "  +7.22%  "
Set-PriceText "D51" "74.57"
